$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1), same bold "header" style as the existing ones ---
$ws.Range("D1").Value = "Email ID"
$ws.Range("E1").Value = "CITY"
$ws.Range("F1").Value = "PINCODE"
$ws.Range("D1:F1").Font.Bold = $true

# --- Column A: First Name values for the 4 data rows ---
$ws.Range("A2").Value = "test101"
$ws.Range("A3").Value = "test102"
$ws.Range("A4").Value = "test103"
$ws.Range("A5").Value = "test104"

# --- Column B: Last Name (same surname reused for every row) ---
$ws.Range("B2").Value = "Sharma"
$ws.Range("B3").Value = "Sharma"
$ws.Range("B4").Value = "Sharma"
$ws.Range("B5").Value = "Sharma"

# --- Column C: Mobile Number ---
$ws.Range("C2").Value = 784352652
$ws.Range("C3").Value = 784352653
$ws.Range("C4").Value = 784352654
$ws.Range("C5").Value = 784352655

# --- Column D: Email ID, entered + hyperlinked one at a time (mirrors the
#     fill-handle-then-overwrite pattern that produced the stray D3:D5 link) ---
$ws.Range("D2").Value = "t101@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:t101@gmail.com")

$ws.Range("D3").Value = "t101@gmail.com"
$ws.Range("D4").Value = "t101@gmail.com"
$ws.Range("D5").Value = "t101@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D3:D5"), "mailto:t101@gmail.com")

$ws.Range("D3").Value = "t102@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:t102@gmail.com")

$ws.Range("D4").Value = "t103@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:t103@gmail.com")

$ws.Range("D5").Value = "t104@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:t104@gmail.com")

# --- Column E: CITY ---
$ws.Range("E2").Value = "BSR"
$ws.Range("E3").Value = "BSR"
$ws.Range("E4").Value = "BLR"
$ws.Range("E5").Value = "BLR"

# --- Column F: PINCODE ---
$ws.Range("F2").Value = 203001
$ws.Range("F3").Value = 203001
$ws.Range("F4").Value = 201005
$ws.Range("F5").Value = 201005

# --- Column D width ---
$ws.Range("D1").ColumnWidth = 19.1

# --- Final selection, matches where the cursor lands after typing the last row ---
$ws.Range("F6").Select()
